$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 14 (the duplicate "clown noes666" alias row); rows below shift up.
$ws.Rows(14).Delete()

# Update shop candidate id (G4) and the mirrored id now sitting in L15.
$ws.Range("G4").Value = 877
$ws.Range("L15").Value = 877

# Clear the alias text that shifted into C14 (was "Clown_noes666").
$ws.Range("C14").Value = ""

# Swap the featured cosmetic: Count Olaf -> Meow Face.
$ws.Range("D6").Value = "Meow Face`nAnimal Crossing"
$ws.Range("G6").Value = ""
$ws.Range("G8").Value = "https://cdn.discordapp.com/attachments/699111007649398865/1062644744498454618/MeowFaceB.png"
$ws.Range("G9").Value = ""

# Update the oldest-creator info block.
$ws.Range("C12").Value = "darkpulse91"
$ws.Range("D12").Value = 72148581
$ws.Range("E12").Value = "224 days"
